# GlobalMart Inc. Backstory and Use Cases document update:
#  - Retitle the document and rename "Business Problem" -> "Financial Problem"
#  - Append a new "Column Descriptions" section documenting each dataset column
$d = $word.ActiveDocument

# --- Title and heading text updates ---
$d.Content.Find.Execute("GlobalMart Inc. Backstory and Analysis Use Cases", $true, $false, $false, $false, $false, $true, 1, $false, "GlobalMart Inc. Backstory, Financial Problem, and Use Case Analysis", 2) | Out-Null
$d.Content.Find.Execute("Business Problem", $true, $false, $false, $false, $false, $true, 1, $false, "Financial Problem", 2) | Out-Null

# --- Append new "Column Descriptions" section with one Heading3 + description per column ---
$columns = @(
    @("Transaction_ID", "A unique identifier for each transaction."),
    @("Date", "The date when the transaction occurred."),
    @("Region", "The geographical region where the transaction took place."),
    @("Product", "The category of the product sold."),
    @("Customer_Age", "The age of the customer who made the purchase."),
    @("Customer_Gender", "The gender of the customer (e.g., Male, Female, Non-Binary)."),
    @("Customer_Segment", "The customer segment (e.g., Retail, Wholesale, Enterprise)."),
    @("Sales_Channel", "The channel through which the sale occurred (e.g., Online, Offline)."),
    @("Units_Sold", "The number of units of the product sold in the transaction."),
    @("Unit_Price", "The price per unit of the product sold."),
    @("Discount", "The discount applied to the unit price of the product."),
    @("Tax", "The tax applied to the transaction."),
    @("Revenue", "The total revenue generated from the transaction."),
    @("Expenses", "The total expenses incurred for the transaction."),
    @("Profit", "The net profit generated from the transaction."),
    @("Profit_Margin", "The percentage of profit relative to revenue."),
    @("Year", "The year in which the transaction occurred."),
    @("Month", "The month in which the transaction occurred."),
    @("Day_of_Week", "The day of the week when the transaction occurred."),
    @("Operational_Inefficiency", "A flag indicating whether the transaction had high operational inefficiency.")
)

$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Heading2"
$p.Range.Text = "Column Descriptions"

foreach ($col in $columns) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Style = "Heading3"
    $p.Range.Text = $col[0]

    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Style = "Normal"
    $p.Range.Text = $col[1]
}

# --- Final trailing empty paragraph (matches trailing <w:p/> in the source) ---
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Normal"
